$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Refresh the time_taken timestamps on the "data" sheet (column F, rows 2-39)
$newTimes = @{
    2  = "2021-10-05 14:34:55.129227"
    3  = "2021-10-05 14:34:55.129235"
    4  = "2021-10-05 14:34:55.129238"
    5  = "2021-10-05 14:34:55.129240"
    6  = "2021-10-05 14:34:55.129243"
    7  = "2021-10-05 14:34:55.129246"
    8  = "2021-10-05 14:34:55.129249"
    9  = "2021-10-05 14:34:55.129251"
    10 = "2021-10-05 14:34:55.129254"
    11 = "2021-10-05 14:34:55.129256"
    12 = "2021-10-05 14:34:55.129259"
    13 = "2021-10-05 14:34:55.129261"
    14 = "2021-10-05 14:34:55.129264"
    15 = "2021-10-05 14:34:55.129266"
    16 = "2021-10-05 14:34:55.129269"
    17 = "2021-10-05 14:34:55.129271"
    18 = "2021-10-05 14:34:55.129274"
    19 = "2021-10-05 14:34:55.129277"
    20 = "2021-10-05 14:34:55.129279"
    21 = "2021-10-05 14:34:55.129282"
    22 = "2021-10-05 14:34:55.129284"
    23 = "2021-10-05 14:34:55.129287"
    24 = "2021-10-05 14:34:55.129289"
    25 = "2021-10-05 14:34:55.129292"
    26 = "2021-10-05 14:34:55.129295"
    27 = "2021-10-05 14:34:55.129297"
    28 = "2021-10-05 14:34:55.129300"
    29 = "2021-10-05 14:34:55.129302"
    30 = "2021-10-05 14:34:55.129304"
    31 = "2021-10-05 14:34:55.129307"
    32 = "2021-10-05 14:34:55.129310"
    33 = "2021-10-05 14:34:55.129312"
    34 = "2021-10-05 14:34:55.129315"
    35 = "2021-10-05 14:34:55.129318"
    36 = "2021-10-05 14:34:55.129320"
    37 = "2021-10-05 14:34:55.129323"
    38 = "2021-10-05 14:34:55.129325"
    39 = "2021-10-05 14:34:55.129328"
}

foreach ($r in $newTimes.Keys) {
    $ws.Cells.Item($r, 6).Value = $newTimes[$r]
}

# Add the new "metadata" sheet right after "data"
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) - text labels
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Copy the bold/bordered header style from the "data" sheet so the new
# header cells match exactly (same style index, no new styles created).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row (row 2)
$meta.Cells.Item(2, 1).Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats, reuse index-column style

$meta.Cells.Item(2, 2).Value = "Mosaic skin disorders"
$meta.Cells.Item(2, 3).Value = 3472

# "1.0" must be stored as text, not converted to the number 1.
$meta.Range("D2").NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.0"
$meta.Range("D2").ClearFormats()

$meta.Cells.Item(2, 5).Value = "2021-02-18T09:45:43.501255Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:34:55.125495"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3472/?format=json"

# Reset selection on metadata sheet to A1 as in source
$meta.Range("A1").Select()
$excel.Application.CutCopyMode = $false

# Keep "data" as the active sheet/tab, matching the unchanged bookViews
# in the source workbook (only the <sheets> list gained a new entry).
$ws.Activate()
